$d = $word.ActiveDocument

# 1. Remove the stale "_GoBack" bookmark that currently sits at the end of the
#    paragraph ending in "...selected for the network to process."
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# 2. Find the paragraph that ends the "Bollinger Bands..." sentence and add a
#    brand-new paragraph after it with the homoscedasticity note.
$finder = $d.Content.Find
$finder.Execute("the more oversold the market.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$insertRange = $finder.Parent
$insertRange.Collapse(0)
$insertRange.InsertParagraphAfter()
$insertRange.Collapse(0)

$newPara = $d.Paragraphs($d.Paragraphs.Count - 1)
$noteText = "A time series is homoscedastic when the variance of the error term (random variable, noise) is constant. Heteroscedastic models don’t assume this."
$newPara.Range.Text = $noteText
$newPara.Format.Alignment = 3

# 3. Re-add the "_GoBack" bookmark as a zero-length mark right after the new
#    text (before the paragraph mark), matching where Word leaves it after
#    edits. A bookmark can't be created directly as a zero-length range that
#    sits exactly on the last character of a paragraph, so we bookmark the
#    final character (non-empty range), remove that character from inside the
#    bookmark (which correctly collapses it in place) and retype it.
$lastCharStart = $newPara.Range.End - 2
$lastCharEnd = $newPara.Range.End - 1
$lastCharRange = $d.Range($lastCharStart, $lastCharEnd)
$lastChar = $lastCharRange.Text
$d.Bookmarks.Add("_GoBack", $lastCharRange)
$bmRange = $d.Bookmarks("_GoBack").Range
$bmRange.Text = ""
$d.Range($lastCharStart, $lastCharStart).InsertAfter($lastChar)
